$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price (D) column cells keep their original text representation
# (values like "1.012" or "0.4717" must not be reinterpreted as numbers).
$priceCells = @("D2","D3","D4","D5","D7","D8","D9","D11","D12","D13","D14","D15","D16","D17","D18","D21","D22","D24","D25","D26","D27","D28","D29","D30","D31","D32","D33","D34","D35","D38","D39","D40","D41","D42","D43","D44","D45","D46","D47","D49","D50","D51")
foreach ($addr in $priceCells) { $ws.Range($addr).NumberFormat = "@" }

$ws.Range("D2").Value = "27.238.87"
$ws.Range("E2").Value = "  +1.19%  "
$ws.Range("D3").Value = "1.835.80"
$ws.Range("E3").Value = "  +1.01%  "
$ws.Range("D4").Value = "1.012"
$ws.Range("E4").Value = "  +1.09%  "
$ws.Range("D5").Value = "314.17"
$ws.Range("E5").Value = "  +1.43%  "
$ws.Range("E6").Value = "  +1.12%  "
$ws.Range("D7").Value = "0.4717"
$ws.Range("E7").Value = "  +0.82%  "
$ws.Range("D8").Value = "0.3692"
$ws.Range("E8").Value = "  -0.33%  "
$ws.Range("D9").Value = "0.07430"
$ws.Range("E9").Value = "  +0.74%  "
$ws.Range("E10").Value = "  +1.29%  "
$ws.Range("D11").Value = "20.47"
$ws.Range("E11").Value = "  +0.07%  "
$ws.Range("D12").Value = "1.829.62"
$ws.Range("E12").Value = "  -1.75%  "
$ws.Range("D13").Value = "0.07344"
$ws.Range("E13").Value = "  +3.71%  "
$ws.Range("D14").Value = "5.482"
$ws.Range("E14").Value = "  +2.09%  "
$ws.Range("D15").Value = "93.15"
$ws.Range("E15").Value = "  +0.40%  "
$ws.Range("D16").Value = "6.580"
$ws.Range("E16").Value = "  +1.06%  "
$ws.Range("D17").Value = "1.014"
$ws.Range("E17").Value = "  +1.24%  "
$ws.Range("D18").Value = "0.000008821"
$ws.Range("E18").Value = "  +1.12%  "
$ws.Range("E19").Value = "  +1.10%  "
$ws.Range("E20").Value = "  +0.38%  "
$ws.Range("D21").Value = "27.256.15"
$ws.Range("E21").Value = "  +1.10%  "
$ws.Range("D22").Value = "5.314"
$ws.Range("E22").Value = "  -0.37%  "
$ws.Range("E23").Value = "  +1.23%  "
$ws.Range("D24").Value = "2.061.54"
$ws.Range("E24").Value = "  +0.00%  "
$ws.Range("D25").Value = "1.904"
$ws.Range("E25").Value = "  +0.04%  "
$ws.Range("D26").Value = "153.15"
$ws.Range("E26").Value = "  +1.05%  "
$ws.Range("D27").Value = "18.65"
$ws.Range("E27").Value = "  +1.39%  "
$ws.Range("D28").Value = "2.173"
$ws.Range("E28").Value = "  -2.11%  "
$ws.Range("D29").Value = "5.279"
$ws.Range("E29").Value = "  -0.91%  "
$ws.Range("D30").Value = "117.95"
$ws.Range("E30").Value = "  +2.04%  "
$ws.Range("D31").Value = "0.08929"
$ws.Range("E31").Value = "  +0.10%  "
$ws.Range("D32").Value = "0.7612"
$ws.Range("E32").Value = "  -0.96%  "
$ws.Range("D33").Value = "1.174"
$ws.Range("E33").Value = "  +0.82%  "
$ws.Range("D34").Value = "4.554"
$ws.Range("E34").Value = "  +1.26%  "
$ws.Range("D35").Value = "2.948"
$ws.Range("E35").Value = "  +1.25%  "
$ws.Range("E36").Value = "  +1.23%  "
$ws.Range("E37").Value = "  +1.73%  "
$ws.Range("D38").Value = "0.05346"
$ws.Range("E38").Value = "  +1.11%  "
$ws.Range("D39").Value = "0.01962"
$ws.Range("E39").Value = "  -0.06%  "
$ws.Range("D40").Value = "2.995"
$ws.Range("E40").Value = "  +1.18%  "
$ws.Range("D41").Value = "7.335"
$ws.Range("E41").Value = "  +0.32%  "
$ws.Range("D42").Value = "2.408"
$ws.Range("E42").Value = "  +1.39%  "
$ws.Range("D43").Value = "0.5352"
$ws.Range("E43").Value = "  -0.04%  "
$ws.Range("D44").Value = "0.1666"
$ws.Range("E44").Value = "  -0.13%  "
$ws.Range("D45").Value = "8.556"
$ws.Range("E45").Value = "  +1.01%  "
$ws.Range("D46").Value = "0.4960"
$ws.Range("E46").Value = "  +0.05%  "
$ws.Range("D47").Value = "10.52"
$ws.Range("E47").Value = "  +0.50%  "
$ws.Range("E48").Value = "  +1.30%  "
$ws.Range("D49").Value = "1.676"
$ws.Range("E49").Value = "  +0.09%  "
$ws.Range("D50").Value = "103.90"
$ws.Range("E50").Value = "  +1.02%  "
$ws.Range("D51").Value = "0.06319"
$ws.Range("E51").Value = "  +0.34%  "
